# Re-generated quadratic/linear problem data (alpha_zero stationary generator).
# Writes new Expression / Function_Evaluation / value cells across the
# "Restricciones_del_follower", "Punto_modificado", "Vector_bf" and
# "Vector_BF" sheets, matching a fresh run of the experiment generator.
#
# Cells that must hold a numeric-looking value as literal TEXT (the workbook
# stores every data value as a shared string, never a native number) are
# written with a leading apostrophe so Excel keeps them as text instead of
# silently recasting them as numbers; the cell style is then reset back to
# "Normal" so no stray quote-prefix formatting lingers on the cell.
function Set-TextValue {
    param($Range, [string]$Text)

    if ($Text -match '^-?[0-9]+(\.[0-9]+)?$') {
        $Range.Value = "'" + $Text
        $Range.Style = "Normal"
    } else {
        $Range.Value = $Text
    }
}

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_follower ---------------------------------------
$ws = $wb.Worksheets.Item(3)   # "Restricciones_del_follower"

Set-TextValue $ws.Cells.Item(2, 1) "-0.25 - x + y_1 + y_2"
Set-TextValue $ws.Cells.Item(2, 2) "0.25"
Set-TextValue $ws.Cells.Item(2, 4) "0.34"
Set-TextValue $ws.Cells.Item(2, 5) "2.3000000000000003"
Set-TextValue $ws.Cells.Item(2, 6) "4.2"

Set-TextValue $ws.Cells.Item(3, 1) "1.7000000000000002 - y_1"
Set-TextValue $ws.Cells.Item(3, 2) "-1.7000000000000002"
Set-TextValue $ws.Cells.Item(3, 4) "0.14"
Set-TextValue $ws.Cells.Item(3, 5) "4.8"
Set-TextValue $ws.Cells.Item(3, 6) "0"

Set-TextValue $ws.Cells.Item(4, 1) "-3.8 - y_2"
Set-TextValue $ws.Cells.Item(4, 2) "-3.8"
Set-TextValue $ws.Cells.Item(4, 4) "0.38"
Set-TextValue $ws.Cells.Item(4, 5) "1.2"
Set-TextValue $ws.Cells.Item(4, 6) "0"

# --- Punto_modificado ---------------------------------------------------
$ws = $wb.Worksheets.Item(4)   # "Punto_modificado"

Set-TextValue $ws.Cells.Item(2, 1) "5.25"
Set-TextValue $ws.Cells.Item(2, 2) "1.7000000000000002"
Set-TextValue $ws.Cells.Item(2, 3) "3.8"

# --- Vector_bf --------------------------------------------------------
# NOTE: worksheet name lookup via Worksheets.Item("...") is case
# insensitive, and this workbook has two sheets whose names differ only
# by case ("Vector_bf" vs "Vector_BF"). Address them by their fixed
# sheet index instead so each edit lands on the intended sheet.
$ws = $wb.Worksheets.Item(5)   # "Vector_bf"

Set-TextValue $ws.Cells.Item(2, 1) "-1.9000000000000001"
Set-TextValue $ws.Cells.Item(3, 1) "-0.96"

# --- Vector_BF ----------------------------------------------------------
$ws = $wb.Worksheets.Item(6)   # "Vector_BF"

Set-TextValue $ws.Cells.Item(2, 1) "2.3000000000000003"
Set-TextValue $ws.Cells.Item(3, 1) "2.7999999999999994"
Set-TextValue $ws.Cells.Item(4, 1) "-2.9000000000000004"

Write-Output "applied generator refresh"
